$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.476.78"
$ws.Range("E2").Value = "  +4.67%  "
$ws.Range("D3").Value = "1.722.88"
$ws.Range("E3").Value = "  +4.01%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.005"
$ws.Range("E4").Value = "  +0.31%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "225.19"
$ws.Range("E5").Value = "  +2.63%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5340"
$ws.Range("E6").Value = "  +1.85%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.005"
$ws.Range("E7").Value = "  +0.23%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2659"
$ws.Range("E8").Value = "  +0.80%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06590"
$ws.Range("E9").Value = "  +3.79%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.44"
$ws.Range("E10").Value = "  +4.46%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07675"
$ws.Range("E11").Value = "  +0.17%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.598"
$ws.Range("E12").Value = "  -0.51%  "
$ws.Range("D13").Value = "1.728.23"
$ws.Range("E13").Value = "  +4.40%  "
$ws.Range("D14").Value = "1.962.52"
$ws.Range("E14").Value = "  +4.29%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5793"
$ws.Range("E15").Value = "  +2.90%  "
$ws.Range("D16").Value = "0.0₅8268"
$ws.Range("E16").Value = "  +1.17%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "67.76"
$ws.Range("E17").Value = "  +3.43%  "
$ws.Range("D18").Value = "27.475.96"
$ws.Range("E18").Value = "  +4.94%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "218.12"
$ws.Range("E19").Value = "  +11.61%  "
$ws.Range("E20").Value = "  +0.15%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.734"
$ws.Range("E21").Value = "  +2.33%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.56"
$ws.Range("E22").Value = "  -0.04%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.998"
$ws.Range("E23").Value = "  +0.47%  "
$ws.Range("E24").Value = "  +0.19%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "143.53"
$ws.Range("E25").Value = "  -1.18%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.733"
$ws.Range("E26").Value = "  +11.14%  "
$ws.Range("E27").Value = "  +3.51%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.318"
$ws.Range("E28").Value = "  +1.05%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "16.46"
$ws.Range("E29").Value = "  +2.99%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05436"
$ws.Range("E30").Value = "  -1.41%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.299"
$ws.Range("E31").Value = "  +2.07%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.545"
$ws.Range("E32").Value = "  +1.63%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.427"
$ws.Range("E33").Value = "  +2.21%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.650"
$ws.Range("E34").Value = "  +5.35%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.877"
$ws.Range("E35").Value = "  +3.07%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9564"
$ws.Range("E36").Value = "  +0.87%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.429"
$ws.Range("E37").Value = "  +0.80%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5920"
$ws.Range("E38").Value = "  +4.69%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01641"
$ws.Range("E39").Value = "  +3.55%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.907"
$ws.Range("E40").Value = "  +2.89%  "
$ws.Range("D41").Value = "1.047.73"
$ws.Range("E41").Value = "  +1.73%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8449"
$ws.Range("E42").Value = "  +2.92%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.005"
$ws.Range("E43").Value = "  +0.22%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "101.17"
$ws.Range("E44").Value = "  +0.43%  "
$ws.Range("D45").Value = "1.867.94"
$ws.Range("E45").Value = "  +4.23%  "
$ws.Range("D46").Value = "0.0₈114"
$ws.Range("E46").Value = "  +5.41%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "58.52"
$ws.Range("E47").Value = "  +1.38%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4510"
$ws.Range("E48").Value = "  +3.96%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.187"
$ws.Range("E49").Value = "  +3.31%  "
$ws.Range("E50").Value = "  +0.38%  "
$ws.Range("B51").Value = "XinFinNetwork"
$ws.Range("C51").Value = "https://coinranking.com/coin/77jGXSqWJ1ofG+xinfinnetwork-xdc"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06521"
$ws.Range("E51").Value = "  +12.38%  "
